# Word COM-interop edit script
# Applies the commit "n complete flakes, platform tpyes and shapes"
$d = $word.ActiveDocument

# Locate paragraphs by a stable leading-text fragment (index-fragile otherwise)
function Find-ParaByPrefix([object]$doc, [string]$prefix) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

# 1. Front-matter date: "21 November, 2016" -> "06 December, 2016"
$pDate = Find-ParaByPrefix $d 'te: "21 November, 2016"'
$pDate.Range.Text = 'te: "06 December, 2016"'

# 2. Introduction paragraph: drop the stray trailing " git"
$pIntro = Find-ParaByPrefix $d 'The Guanyindong site, located'
$pIntro.Range.Text = 'The Guanyindong site, located in Guanyindong village, Qianxi County of Guizhou Province (26°51′26″N, 105°58′7″E) at an elevation of 1464 m a.s.l., is a limestone cave site extending from east to west it was discovered by a team organized by the institute of Vertebrate Paleontology and Paleoanthropolgy(IVPP),Chinese Academy of Sciences in 1964. Several excavations were conducted in 1965, 1972 and 1973, yeilding A total of 176 cores, 1292 flakes, 1101 retouched pieces and 802 pieces of debris were identified.'

# 3. Flakes paragraph: "We found xx complete flakes." -> "We found 196 complete unretouched flakes."
$pFlakes = Find-ParaByPrefix $d 'We found xx complete flakes.'
$pFlakes.Range.Text = 'We found 196 complete unretouched flakes. The average maximum length of the flakes is 62.6 mm, the average thickness is``mm. There are xxx flakes or flake breaks that have distinguishable platform, that can be divided into cortex(%), plain(%), facet(%),dihederal(%) and focus(%). The shapes of platform include triangle (%), fusiform (%), quadrangle (%) and CDG (%). Most of flakes dorsal side is partially covered with cortex (%). And the average scar number is xxx. The directions of these scars are recorded. Among them, the scars with the same directions of flake are dominated (%), following with opposite direction(%). We also found a number of centripetal scars (%).'

# 4. Retouched pieces paragraph: extend with platform-type / shape sentences
$pRetouched = Find-ParaByPrefix $d 'A total of 1101 retouched pieces were found'
$pRetouched.Range.Text = 'A total of 1101 retouched pieces were found, accounting for 48.5% of lithic assemblage. The average max dimension is xxx. % retouched pieces are made on flakes (%) and flake breaks (%), others are made on either chunks or pebbles. Side scrapers dominate the sub-division of retouched pieces (%), followed by denticulates and borers. Convex edge constitutes the largest proportion of the edge shapes of side scrapers(%). Looking at the location of retouch and the size of the retouched flakes can provide us further insight into retouching behaviours. Most of tools have more than one retouched edges. We introduced two concepts “Zone Index” and “Geometric Index of Unifacial Reduction(GIUR)” to estimate the invasion and intensity of retouching. From our observation… We also measured the angle of each retouched edge. For notch pieces(n=91), we found that most notches only have one notch end on each retouched piece and the average depth and length is xx and xx. The location of retouching is mainly on one side which defined as longer side of the piece.'

# 5. New paragraph after the "Levallois" heading describing Levallois-like artifacts
$pLevalloisHeading = Find-ParaByPrefix $d 'Levallois'
$newPara = $pLevalloisHeading.Range.InsertParagraphAfter()
$insertedPara = $pLevalloisHeading.Next()
$insertedPara.Style = "First Paragraph"
$insertedPara.Range.Text = 'We distinguished 70 stone artifacts that are Levallois like including 11 Levallois cores, 22 flakes, 4 points and 33 tools made on levallois flakes. The average dimension of levallois products is xxx which is smaller (or larger) than ordinary products. The platform shapes of levallois flakes are various ranging from triangle, quadrangle, fusiformis to chapeau de gendarme(CDG). For flakes, we measured the thickness at 25%, 50%, 70% max dimension and compared them with ordinary flakes found that levallois flakes are relatively more flat. This found is consistent with the theory that ... The scar number is also relatively more than ordinary flakes, the direction of which is mostly centripetal.'

Write-Output "edits applied"
